$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item(1)          # "总计" sheet
$q3 = $wb.Worksheets.Item("2022-Q3")     # sheet currently holding the (now old) 2022-Q3 fund data

# ------------------------------------------------------------------
# 1) Rename the current "2022-Q3" sheet to "2022-Q4" - it keeps its sheetId
#    and will be overwritten below with the new Q4 fund-holdings table.
# ------------------------------------------------------------------
$q4 = $q3
$q4.Name = "2022-Q4"

# ------------------------------------------------------------------
# 2) Insert a brand new sheet right after it, named "2022-Q3", which will
#    hold the data that used to live on the original "2022-Q3" sheet.
# ------------------------------------------------------------------
$q3new = $wb.Worksheets.Add($null, $q4)
$q3new.Name = "2022-Q3"

# Copy the cell formatting (style) that used to be on the Q3 sheet (before we
# overwrite it with Q4 data below) onto the new Q3 sheet, so it keeps looking
# the same as before the edit.
$q4.Range("B1").Copy()
$q3new.Range("B1:H1").PasteSpecial(-4122)
$q4.Range("A2").Copy()
$q3new.Range("A2:A3").PasteSpecial(-4122)

# Match the page margins the original "2022-Q3" sheet used to have.
$q3new.PageSetup.LeftMargin = 50.4
$q3new.PageSetup.RightMargin = 50.4
$q3new.PageSetup.TopMargin = 54
$q3new.PageSetup.BottomMargin = 54
$q3new.PageSetup.HeaderMargin = 21.599999999999998
$q3new.PageSetup.FooterMargin = 21.599999999999998

# ------------------------------------------------------------------
# 3) Re-populate the new "2022-Q3" sheet with the original fund-holdings data
# ------------------------------------------------------------------
$q3new.Range("B1").Value = "基金代码"
$q3new.Range("C1").Value = "基金名称"
$q3new.Range("D1").Value = "基金规模"
$q3new.Range("E1").Value = "股票总仓位"
$q3new.Range("F1").Value = "仓位占比"
$q3new.Range("G1").Value = "持有市值(亿元)"
$q3new.Range("H1").Value = "仓位排名"

$q3new.Range("A2").Value = 0
$q3new.Range("B2").Value = "'011761"
$q3new.Range("C2").Value = "平安鑫瑞混合A"
$q3new.Range("D2").Value = "'0.59"
$q3new.Range("E2").Value = "'23.80"
$q3new.Range("F2").Value = "'1.01"
$q3new.Range("G2").Value = "'0.0060"
$q3new.Range("H2").Value = 7

$q3new.Range("A3").Value = 1
$q3new.Range("B3").Value = "'011762"
$q3new.Range("C3").Value = "平安鑫瑞混合C"
$q3new.Range("D3").Value = "'0.28"
$q3new.Range("E3").Value = "'23.80"
$q3new.Range("F3").Value = "'1.01"
$q3new.Range("G3").Value = "'0.0028"
$q3new.Range("H3").Value = 7

# ------------------------------------------------------------------
# 4) Overwrite the renamed "2022-Q4" sheet with the new Q4 fund-holdings data
# ------------------------------------------------------------------
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "'000968"
$q4.Range("C2").Value = "广发中证养老产业指数A"
$q4.Range("D2").Value = "'12.03"
$q4.Range("E2").Value = "'94.66"
$q4.Range("F2").Value = "'1.59"
$q4.Range("G2").Value = "'0.1913"
$q4.Range("H2").Value = 1

$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "'159855"
$q4.Range("C3").Value = "银华中证影视主题ETF"
$q4.Range("D3").Value = "'1.01"
$q4.Range("E3").Value = "'97.80"
$q4.Range("F3").Value = "'4.06"
$q4.Range("G3").Value = "'0.0410"
$q4.Range("H3").Value = 8

$q4.Range("A4").Value = 2
$q4.Range("B4").Value = "'516620"
$q4.Range("C4").Value = "国泰中证影视主题ETF"
$q4.Range("D4").Value = "'0.71"
$q4.Range("E4").Value = "'98.01"
$q4.Range("F4").Value = "'4.17"
$q4.Range("G4").Value = "'0.0296"
$q4.Range("H4").Value = 8

$q4.Range("A5").Value = 3
$q4.Range("B5").Value = "'002982"
$q4.Range("C5").Value = "广发中证养老产业指数C"
$q4.Range("D5").Value = "'0.97"
$q4.Range("E5").Value = "'94.66"
$q4.Range("F5").Value = "'1.59"
$q4.Range("G5").Value = "'0.0154"
$q4.Range("H5").Value = 1

$q4.Range("A6").Value = 4
$q4.Range("B6").Value = "'516560"
$q4.Range("C6").Value = "华宝养老ETF"
$q4.Range("D6").Value = "'0.79"
$q4.Range("E6").Value = "'98.20"
$q4.Range("F6").Value = "'1.68"
$q4.Range("G6").Value = "'0.0133"
$q4.Range("H6").Value = 1

$q4.Range("A7").Value = 5
$q4.Range("B7").Value = "'165524"
$q4.Range("C7").Value = "信诚中证智能家居指数（LOF）A"
$q4.Range("D7").Value = "'0.37"
$q4.Range("E7").Value = "'91.28"
$q4.Range("F7").Value = "'1.29"
$q4.Range("G7").Value = "'0.0048"
$q4.Range("H7").Value = 1

$q4.Range("A8").Value = 6
$q4.Range("B8").Value = "'013084"
$q4.Range("C8").Value = "信诚中证智能家居指数（LOF）C"
$q4.Range("D8").Value = "'0.14"
$q4.Range("E8").Value = "'91.28"
$q4.Range("F8").Value = "'1.29"
$q4.Range("G8").Value = "'0.0018"
$q4.Range("H8").Value = 1

# Re-apply the "2022-Q4" header/row-label formatting to match the "总计"
# sheet's look (same style used for the new Q4 data as in the source workbook).
$total.Range("B1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$total.Range("A2").Copy()
$q4.Range("A2:A8").PasteSpecial(-4122)

# Match the page margins that "总计" uses, since the new Q4 sheet mirrors it.
$q4.PageSetup.LeftMargin = 54
$q4.PageSetup.RightMargin = 54
$q4.PageSetup.TopMargin = 72
$q4.PageSetup.BottomMargin = 72
$q4.PageSetup.HeaderMargin = 36
$q4.PageSetup.FooterMargin = 36

# ------------------------------------------------------------------
# 5) Update the "总计" (totals) sheet: move the old Q3 total down to row 3,
#    and put the new Q4 total into row 2.
# ------------------------------------------------------------------
$total.Range("A2").Copy($total.Range("A3"))
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q3"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.01

$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 7
$total.Range("D2").Value = 0.3

# Keep the "2022-Q3" sheet as the visually selected tab, matching the source.
$q3new.Activate()
